# Fruta / hortaliza, semanal
# Insert 3 new weekly rows (Feria Lagunitas de Puerto Montt - Pera) above the
# current row 139, pushing the existing rows 139-152 down to 142-155.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("139:141").Insert()

# Row 139: new Forelle entry
$ws.Cells.Item(139, 1).Value = 4
$ws.Cells.Item(139, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(139, 3).Value = "Los Lagos"
$ws.Cells.Item(139, 4).Value = 44491
$ws.Cells.Item(139, 5).Value = 10
$ws.Cells.Item(139, 6).Value = "Fruta"
$ws.Cells.Item(139, 7).Value = 100104
$ws.Cells.Item(139, 8).Value = "Frutos de pepita"
$ws.Cells.Item(139, 9).Value = 100104005
$ws.Cells.Item(139, 10).Value = "Pera"
$ws.Cells.Item(139, 11).Value = "Forelle"
$ws.Cells.Item(139, 12).Value = "Primera"
$ws.Cells.Item(139, 13).Value = 400
$ws.Cells.Item(139, 14).Value = 13000
$ws.Cells.Item(139, 15).Value = 13500
$ws.Cells.Item(139, 16).Value = 13250
$ws.Cells.Item(139, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(139, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(139, 19).Value = 883
$ws.Cells.Item(139, 20).Value = 15

# Row 140: new Packham's Triumph / Primera entry
$ws.Cells.Item(140, 1).Value = 4
$ws.Cells.Item(140, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(140, 3).Value = "Los Lagos"
$ws.Cells.Item(140, 4).Value = 44491
$ws.Cells.Item(140, 5).Value = 10
$ws.Cells.Item(140, 6).Value = "Fruta"
$ws.Cells.Item(140, 7).Value = 100104
$ws.Cells.Item(140, 8).Value = "Frutos de pepita"
$ws.Cells.Item(140, 9).Value = 100104005
$ws.Cells.Item(140, 10).Value = "Pera"
$ws.Cells.Item(140, 11).Value = "Packham's Triumph"
$ws.Cells.Item(140, 12).Value = "Primera"
$ws.Cells.Item(140, 13).Value = 500
$ws.Cells.Item(140, 14).Value = 15000
$ws.Cells.Item(140, 15).Value = 16000
$ws.Cells.Item(140, 16).Value = 15500
$ws.Cells.Item(140, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(140, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(140, 19).Value = 1033
$ws.Cells.Item(140, 20).Value = 15

# Row 141: new Packham's Triumph / Segunda entry
$ws.Cells.Item(141, 1).Value = 4
$ws.Cells.Item(141, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(141, 3).Value = "Los Lagos"
$ws.Cells.Item(141, 4).Value = 44491
$ws.Cells.Item(141, 5).Value = 10
$ws.Cells.Item(141, 6).Value = "Fruta"
$ws.Cells.Item(141, 7).Value = 100104
$ws.Cells.Item(141, 8).Value = "Frutos de pepita"
$ws.Cells.Item(141, 9).Value = 100104005
$ws.Cells.Item(141, 10).Value = "Pera"
$ws.Cells.Item(141, 11).Value = "Packham's Triumph"
$ws.Cells.Item(141, 12).Value = "Segunda"
$ws.Cells.Item(141, 13).Value = 200
$ws.Cells.Item(141, 14).Value = 13000
$ws.Cells.Item(141, 15).Value = 13000
$ws.Cells.Item(141, 16).Value = 13000
$ws.Cells.Item(141, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(141, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(141, 19).Value = 867
$ws.Cells.Item(141, 20).Value = 15
